# Atualização de bases das ligas, do dia: 24-02-2024 às 23:13
#
# The upstream "Scotland Premiership" odds feed re-pulled match ids
# 6844776/6845513 (24-02-24 Kilmarnock/Hibernian-Aberdeen fixtures),
# 6844802/6844803 (Hearts/Ross County vs Kilmarnock/Dundee) and
# 7667636/7667637 (Aberdeen vs Motherwell fixtures), which swapped the
# two rows' worth of odds data (everything except the running index in
# column A), plus a small odds correction on row 278 (id 276).
#
# This script writes each changed cell's new literal value directly
# (values taken from the authoritative re-pull), cell by cell, matching
# exactly how Excel would record manual/paste edits against the
# existing worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 203 (id 201) picks up what used to be row 204's data ---
$ws.Range("B203").Value = 6845513
$ws.Range("F203").Value = "Hibernian"
$ws.Range("G203").Value = "Aberdeen"
$ws.Range("K203").Value = 2
$ws.Range("L203").Value = 3.4
$ws.Range("M203").Value = 3.3
$ws.Range("N203").Value = 2.3
$ws.Range("O203").Value = 3.5
$ws.Range("P203").Value = 3
$ws.Range("Q203").Value = -0.25
$ws.Range("R203").Value = 2.025
$ws.Range("S203").Value = 1.825
$ws.Range("T203").Value = 2.5
$ws.Range("U203").Value = 1.95
$ws.Range("V203").Value = 1.9
$ws.Range("W203").Value = 1.3
$ws.Range("Z203").Value = 1.025
$ws.Range("AA203").Value = -1
$ws.Range("AC203").Value = 0.8999999999999999

# --- Row 204 (id 202) picks up what used to be row 203's data ---
$ws.Range("B204").Value = 6844776
$ws.Range("F204").Value = "Rangers"
$ws.Range("G204").Value = "St Mirren"
$ws.Range("K204").Value = 1.25
$ws.Range("L204").Value = 5.5
$ws.Range("M204").Value = 8.5
$ws.Range("N204").Value = 1.166
$ws.Range("O204").Value = 7
$ws.Range("P204").Value = 17
$ws.Range("Q204").Value = -2
$ws.Range("R204").Value = 1.875
$ws.Range("S204").Value = 1.975
$ws.Range("T204").Value = 3
$ws.Range("U204").Value = 1.825
$ws.Range("V204").Value = 2.025
$ws.Range("W204").Value = 0.1659999999999999
$ws.Range("Z204").Value = 0
$ws.Range("AA204").Value = -0
$ws.Range("AC204").Value = 1.025

# --- Row 234 (id 232) picks up what used to be row 235's data ---
$ws.Range("B234").Value = 6844802
$ws.Range("F234").Value = "Hearts"
$ws.Range("G234").Value = "Ross County"
$ws.Range("K234").Value = 1.666
$ws.Range("L234").Value = 3.75
$ws.Range("M234").Value = 5
$ws.Range("N234").Value = 1.615
$ws.Range("O234").Value = 4
$ws.Range("P234").Value = 5.25
$ws.Range("R234").Value = 1.8
$ws.Range("S234").Value = 2.05
$ws.Range("U234").Value = 1.95
$ws.Range("V234").Value = 1.9
$ws.Range("X234").Value = 3
$ws.Range("AA234").Value = 1.05
$ws.Range("AB234").Value = 0.95

# --- Row 235 (id 233) picks up what used to be row 234's data ---
$ws.Range("B235").Value = 6844803
$ws.Range("F235").Value = "Kilmarnock"
$ws.Range("G235").Value = "Dundee"
$ws.Range("K235").Value = 1.909
$ws.Range("L235").Value = 3.5
$ws.Range("M235").Value = 3.9
$ws.Range("N235").Value = 1.727
$ws.Range("O235").Value = 3.75
$ws.Range("P235").Value = 4.75
$ws.Range("R235").Value = 1.925
$ws.Range("S235").Value = 1.925
$ws.Range("U235").Value = 1.975
$ws.Range("V235").Value = 1.875
$ws.Range("X235").Value = 2.75
$ws.Range("AA235").Value = 0.925
$ws.Range("AB235").Value = 0.9750000000000001

# --- Row 265 (id 263) picks up what used to be row 266's data ---
$ws.Range("B265").Value = 7667637
$ws.Range("F265").Value = "Aberdeen"
$ws.Range("G265").Value = "Motherwell"
$ws.Range("I265").Value = 3
$ws.Range("J265").Value = "D"
$ws.Range("K265").Value = 1.85
$ws.Range("L265").Value = 3.5
$ws.Range("M265").Value = 4.2
$ws.Range("N265").Value = 1.95
$ws.Range("O265").Value = 3.5
$ws.Range("P265").Value = 3.8
$ws.Range("Q265").Value = -0.5
$ws.Range("T265").Value = 2.5
$ws.Range("U265").Value = 2.05
$ws.Range("V265").Value = 1.8
$ws.Range("W265").Value = -1
$ws.Range("X265").Value = 2.5
$ws.Range("AB265").Value = 1.05
$ws.Range("AC265").Value = -1

# --- Row 266 (id 264) picks up what used to be row 265's data ---
$ws.Range("B266").Value = 7667636
$ws.Range("F266").Value = "Rangers"
$ws.Range("G266").Value = "Ross County"
$ws.Range("I266").Value = 1
$ws.Range("J266").Value = "H"
$ws.Range("K266").Value = 1.1
$ws.Range("L266").Value = 10
$ws.Range("M266").Value = 23
$ws.Range("N266").Value = 1.09
$ws.Range("O266").Value = 12
$ws.Range("P266").Value = 19
$ws.Range("Q266").Value = -2.75
$ws.Range("T266").Value = 3.75
$ws.Range("U266").Value = 1.95
$ws.Range("V266").Value = 1.9
$ws.Range("W266").Value = 0.09000000000000008
$ws.Range("X266").Value = -1
$ws.Range("AB266").Value = 0.475
$ws.Range("AC266").Value = -0.5

# --- Row 278 (id 276): small odds correction, no row swap ---
$ws.Range("R278").Value = 1.95
$ws.Range("S278").Value = 1.9
